$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.047.77"
$ws.Range("E2").Value = "  -3.67%  "
$ws.Range("D3").Value = "3.350.01"
$ws.Range("E3").Value = "  -4.41%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'569.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.52%  "
$ws.Range("D6").Value = "'124.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.29%  "
$ws.Range("D8").Value = "3.353.90"
$ws.Range("E8").Value = "  -4.25%  "
$ws.Range("D9").Value = "'0.474"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.89%  "
$ws.Range("D10").Value = "'7.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.37%  "
$ws.Range("D11").Value = "'0.117"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.27%  "
$ws.Range("D12").Value = "'0.372"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.71%  "
$ws.Range("D13").Value = "3.938.01"
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "3.369.47"
$ws.Range("E15").Value = "  -3.97%  "
$ws.Range("D16").Value = "'0.0000168"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.24%  "
$ws.Range("D17").Value = "62.288.32"
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").Value = "'24.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.79%  "
$ws.Range("D19").Value = "'9.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -10.14%  "
$ws.Range("D20").Value = "'5.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.79%  "
$ws.Range("D21").Value = "'12.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.51%  "
$ws.Range("D22").Value = "'359.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.49%  "
$ws.Range("D23").Value = "'0.548"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.08%  "
$ws.Range("D24").Value = "3.490.89"
$ws.Range("E24").Value = "  -4.26%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "'70.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.08%  "
$ws.Range("D27").Value = "'0.0000104"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -11.22%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").Value = "'6.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.55%  "
$ws.Range("D30").Value = "'1.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.54%  "
$ws.Range("D31").Value = "'7.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.56%  "
$ws.Range("D32").Value = "'2.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.99%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "3.384.45"
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("D35").Value = "'0.146"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.04%  "
$ws.Range("D36").Value = "'22.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.81%  "
$ws.Range("D37").Value = "'5.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").Value = "'165.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "'6.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.02%  "
$ws.Range("D40").Value = "'1.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.34%  "
$ws.Range("D41").Value = "'0.0748"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.57%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "'41.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("D44").Value = "'0.755"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.82%  "
$ws.Range("D45").Value = "'4.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.70%  "
$ws.Range("D46").Value = "'1.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.89%  "
$ws.Range("D47").Value = "'1.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.54%  "
$ws.Range("D48").Value = "'22.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.34%  "
$ws.Range("D49").Value = "'6.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.47%  "
$ws.Range("D50").Value = "2.206.45"
$ws.Range("E50").Value = "  -9.04%  "
$ws.Range("D51").Value = "'0.831"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.50%  "
